$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.894.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.68"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.00%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "330.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3649"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.62%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.73"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "

$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("E10").Value = "  +1.47%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07055"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.42%  "

$ws.Range("E12").Value = "  +0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.068"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.49%  "

$ws.Range("E14").Value = "  +1.49%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.664.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.611"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06549"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.00%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.914"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.58%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.896.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.97%  "

$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.399"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.09"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.847.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.173"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.075"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.56%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.731"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08452"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.18%  "

$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.24"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.140"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.229"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.06%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06023"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02228"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2083"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.86%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.208"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.57%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.57"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.853"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5717"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.62%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("E49").Value = "  +1.19%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06992"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.187"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.31%  "
